$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect text-like numeric strings (e.g. "1.0000", "0.5200") from Excel's
# automatic number coercion by forcing the Price/Volume columns to Text format
# before writing any values.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply the updated cell values from the crypto price refresh
$ws.Range("D2").Value = '26.348.30'
$ws.Range("E2").Value = '  -3.10%  '
$ws.Range("D3").Value = '1.832.48'
$ws.Range("E3").Value = '  -2.63%  '
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '259.32'
$ws.Range("E5").Value = '  -7.84%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '0.5200'
$ws.Range("D8").Value = '0.3217'
$ws.Range("E8").Value = '  -9.00%  '
$ws.Range("D9").Value = '0.06733'
$ws.Range("E9").Value = '  -4.22%  '
$ws.Range("D10").Value = '18.67'
$ws.Range("E10").Value = '  -8.40%  '
$ws.Range("D11").Value = '0.7647'
$ws.Range("E11").Value = '  -7.13%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07674'
$ws.Range("E12").Value = '  -1.83%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.872.61'
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").Value = '88.85'
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").Value = '5.021'
$ws.Range("E15").Value = '  -3.51%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '14.04'
$ws.Range("E17").Value = '  -3.99%  '
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '0.000007890'
$ws.Range("E19").Value = '  -3.53%  '
$ws.Range("D20").Value = '26.409.91'
$ws.Range("E20").Value = '  -2.97%  '
$ws.Range("D21").Value = '2.087.95'
$ws.Range("E21").Value = '  -2.14%  '
$ws.Range("D22").Value = '4.538'
$ws.Range("E22").Value = '  -5.01%  '
$ws.Range("D23").Value = '9.422'
$ws.Range("E23").Value = '  -7.29%  '
$ws.Range("D24").Value = '5.917'
$ws.Range("E24").Value = '  -5.43%  '
$ws.Range("D25").Value = '2.275'
$ws.Range("E25").Value = '  -5.35%  '
$ws.Range("D26").Value = '145.05'
$ws.Range("E26").Value = '  -1.57%  '
$ws.Range("D27").Value = '1.641'
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").Value = '16.90'
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("D29").Value = '111.16'
$ws.Range("E29").Value = '  -2.82%  '
$ws.Range("D30").Value = '4.193'
$ws.Range("E30").Value = '  -5.16%  '
$ws.Range("D31").Value = '4.131'
$ws.Range("E31").Value = '  -6.07%  '
$ws.Range("D32").Value = '0.08741'
$ws.Range("E32").Value = '  -2.53%  '
$ws.Range("E33").Value = '  -2.07%  '
$ws.Range("D34").Value = '1.126'
$ws.Range("E34").Value = '  -5.18%  '
$ws.Range("D35").Value = '2.842'
$ws.Range("E35").Value = '  -2.13%  '
$ws.Range("E36").Value = '  -9.14%  '
$ws.Range("D37").Value = '3.091'
$ws.Range("E37").Value = '  -6.80%  '
$ws.Range("D38").Value = '0.01776'
$ws.Range("E38").Value = '  -5.93%  '
$ws.Range("E39").Value = '  -8.02%  '
$ws.Range("D40").Value = '0.4906'
$ws.Range("E40").Value = '  -7.64%  '
$ws.Range("D41").Value = '112.38'
$ws.Range("E41").Value = '  -3.93%  '
$ws.Range("D42").Value = '0.8895'
$ws.Range("E42").Value = '  -8.66%  '
$ws.Range("D43").Value = '6.125'
$ws.Range("E43").Value = '  -3.21%  '
$ws.Range("D44").Value = '0.9995'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '7.704'
$ws.Range("E45").Value = '  -6.42%  '
$ws.Range("D46").Value = '0.4199'
$ws.Range("E46").Value = '  -9.09%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1254'
$ws.Range("E47").Value = '  -8.62%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05874'
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.037'
$ws.Range("E49").Value = '  -4.84%  '
$ws.Range("D50").Value = '35.31'
$ws.Range("E50").Value = '  -3.88%  '
$ws.Range("D51").Value = '59.33'
$ws.Range("E51").Value = '  -3.98%  '

# Restore the original General format now that the text values are safely stored
$ws.Range("D2:E51").NumberFormat = "General"
